# Update countries & provincias Spain
# Refresh the COVID country table ("Pais" sheet) with the next data pull:
#   - header timestamp (A1) advances from 11 Jul 22:53 to 12 Jul 00:10
#   - several countries received updated totals, which re-sorts a few
#     adjacent rows (Rumania/Barein, Botsuana/Comoras, Monaco/Bahamas)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: updated timestamp
$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 00:10"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3352017
$ws.Range("C4").Value = 58090
$ws.Range("D4").Value = 1487708
$ws.Range("E4").Value = 1726949
$ws.Range("G4").Value = 689
$ws.Range("H4").Value = 137360

# Row 5: Brasil
$ws.Range("B5").Value = 1839850
$ws.Range("C5").Value = 35512
$ws.Range("E5").Value = 554869
$ws.Range("G5").Value = 945
$ws.Range("H5").Value = 71469

# Row 8: Peru
$ws.Range("B8").Value = 322710
$ws.Range("C8").Value = 3064
$ws.Range("D8").Value = 214152
$ws.Range("E8").Value = 96876
$ws.Range("G8").Value = 182
$ws.Range("H8").Value = 11682

# Row 50: now Barein (was Rumania) - totals re-sorted the pair
$ws.Range("A50").Value = "Barein"
$ws.Range("B50").Value = 32470
$ws.Range("C50").Value = 431
$ws.Range("D50").Value = 27828
$ws.Range("E50").Value = 4538
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 104

# Row 51: now Rumania (was Barein)
$ws.Range("A51").Value = "Rumania"
$ws.Range("B51").Value = 32079
$ws.Range("C51").Value = 698
$ws.Range("D51").Value = 21414
$ws.Range("E51").Value = 8794
$ws.Range("G51").Value = 24
$ws.Range("H51").Value = 1871

# Row 71: Costa de Marfil
$ws.Range("B71").Value = 12443
$ws.Range("C71").Value = 391
$ws.Range("D71").Value = 6357
$ws.Range("E71").Value = 6004
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 82

# Row 86: Bulgaria
$ws.Range("B86").Value = 7175
$ws.Range("C86").Value = 211
$ws.Range("D86").Value = 3311
$ws.Range("E86").Value = 3597

# Row 93: Estado de Palestina
$ws.Range("E93").Value = 5362
$ws.Range("G93").Value = 6
$ws.Range("H93").Value = 33

# Row 94: Mauritania
$ws.Range("B94").Value = 5275
$ws.Range("C94").Value = 72
$ws.Range("D94").Value = 2160
$ws.Range("E94").Value = 2968
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 147

# Row 139: Burkina Faso
$ws.Range("B139").Value = 1033
$ws.Range("C139").Value = 13
$ws.Range("D139").Value = 869
$ws.Range("E139").Value = 111

# Row 140: Republica de Chipre
$ws.Range("B140").Value = 1014
$ws.Range("C140").Value = 1
$ws.Range("E140").Value = 156

# Row 165: now Comoras (was Botsuana) - totals re-sorted the pair
$ws.Range("A165").Value = "Comoras"
$ws.Range("B165").Value = 317
$ws.Range("C165").Value = 3
$ws.Range("D165").Value = 296
$ws.Range("E165").Value = 14
$ws.Range("H165").Value = 7

# Row 166: now Botsuana (was Comoras)
$ws.Range("A166").Value = "Botsuana"
$ws.Range("D166").Value = 31
$ws.Range("E166").Value = 282
$ws.Range("H166").Value = 1

# Row 181: now Bahamas (was Monaco) - totals re-sorted the pair
$ws.Range("A181").Value = "Bahamas"
$ws.Range("B181").Value = 111
$ws.Range("C181").Value = 4
$ws.Range("D181").Value = 89
$ws.Range("E181").Value = 11
$ws.Range("H181").Value = 11

# Row 182: now Monaco (was Bahamas)
$ws.Range("A182").Value = "Monaco"
$ws.Range("B182").Value = 109
$ws.Range("D182").Value = 96
$ws.Range("E182").Value = 9
$ws.Range("H182").Value = 4
